$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.288.33"
$ws.Range("E2").Value = "  +2.42%  "

$ws.Range("D3").Value = "1.879.27"
$ws.Range("E3").Value = "  +5.04%  "

$cell = $ws.Range("D4")
$cell.Value = "'1.002"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$cell = $ws.Range("D5")
$cell.Value = "'312.27"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "

$cell = $ws.Range("D6")
$cell.Value = "'1.003"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$cell = $ws.Range("D7")
$cell.Value = "'0.5064"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.42%  "

$cell = $ws.Range("D8")
$cell.Value = "'0.3945"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.66%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.09634"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +4.10%  "

$cell = $ws.Range("D10")
$cell.Value = "'1.146"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'40.91"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +1.12%  "

$cell = $ws.Range("D12")
$cell.Value = "'6.492"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +4.36%  "

$cell = $ws.Range("D13")
$cell.Value = "'21.01"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("D14").Value = "1.879.13"
$ws.Range("E14").Value = "  +4.83%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D15")
$cell.Value = "'7.449"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +4.74%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$cell = $ws.Range("D16")
$cell.Value = "'1.002"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$cell = $ws.Range("D17")
$cell.Value = "'0.00001129"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.46%  "

$cell = $ws.Range("D18")
$cell.Value = "'93.13"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.61%  "

$cell = $ws.Range("D19")
$cell.Value = "'0.06599"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "

$cell = $ws.Range("D20")
$cell.Value = "'17.65"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.02%  "

$cell = $ws.Range("D21")
$cell.Value = "'1.002"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "

$cell = $ws.Range("D22")
$cell.Value = "'6.204"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +5.40%  "

$ws.Range("D23").Value = "28.336.41"
$ws.Range("E23").Value = "  +2.47%  "

$cell = $ws.Range("D24")
$cell.Value = "'11.34"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.40%  "

$cell = $ws.Range("D25")
$cell.Value = "'2.303"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.58%  "

$cell = $ws.Range("D26")
$cell.Value = "'2.574"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +8.43%  "

$ws.Range("D27").Value = "2.096.12"
$ws.Range("E27").Value = "  +4.94%  "

$ws.Range("E28").Value = "  +5.06%  "

$cell = $ws.Range("D29")
$cell.Value = "'159.07"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "

$cell = $ws.Range("D30")
$cell.Value = "'127.66"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.41%  "

$cell = $ws.Range("D31")
$cell.Value = "'0.1072"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("E32").Value = "  +2.70%  "

$cell = $ws.Range("D33")
$cell.Value = "'5.647"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.00%  "

$cell = $ws.Range("D34")
$cell.Value = "'3.625"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$cell = $ws.Range("D35")
$cell.Value = "'9.572"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +9.18%  "

$cell = $ws.Range("D36")
$cell.Value = "'0.06724"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.61%  "

$cell = $ws.Range("D37")
$cell.Value = "'0.02391"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.78%  "

$ws.Range("E38").Value = "  +4.15%  "

$cell = $ws.Range("D39")
$cell.Value = "'11.53"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.24%  "

$cell = $ws.Range("D40")
$cell.Value = "'0.6371"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.80%  "

$cell = $ws.Range("D41")
$cell.Value = "'5.012"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +2.83%  "

$cell = $ws.Range("D42")
$cell.Value = "'1.189"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +4.48%  "

$cell = $ws.Range("D43")
$cell.Value = "'1.002"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

$cell = $ws.Range("D44")
$cell.Value = "'13.46"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.49%  "

$cell = $ws.Range("D45")
$cell.Value = "'0.5992"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.12%  "

$cell = $ws.Range("D46")
$cell.Value = "'3.663"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("E47").Value = "  +0.35%  "

$cell = $ws.Range("D48")
$cell.Value = "'2.009"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +5.27%  "

$cell = $ws.Range("D49")
$cell.Value = "'124.37"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("E50").Value = "  +2.87%  "

$cell = $ws.Range("D51")
$cell.Value = "'0.06865"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.69%  "
